# Fruta / hortaliza, semanal
# A new weekly record was added to the "Tuna" sheet. In the canonical
# OOXML this shows up as a brand-new row 26 (date serial 45002) with all
# of the previously-existing rows 26-80 pushed down by one (to 27-81).
#
# Reproduce that with a real row insert so every existing record keeps
# its values untouched, then populate the newly-opened row 26 with the
# new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 26 - this shifts the
# old rows 26..80 down to 27..81 and bumps the sheet dimension to T81.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly observation.
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value = 45002
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100107
$ws.Range("H26").Value = "Otros"
$ws.Range("I26").Value = 100107011
$ws.Range("J26").Value = "Tuna"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 170
$ws.Range("N26").Value = 18000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 19059
$ws.Range("Q26").Value = "$/caja 16 kilos"
$ws.Range("R26").Value = "Provincia de Los Andes"
$ws.Range("S26").Value = 1191
$ws.Range("T26").Value = 16
